$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.505.16"
$ws.Range("E2").Value = "  -2.84%  "

# Row 3
$ws.Range("D3").Value = "1.657.62"
$ws.Range("E3").Value = "  -4.24%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.38%  "

# Row 6
$ws.Range("E6").Value = "  -2.72%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0618"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0876"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.04%  "

# Row 12
$ws.Range("D12").Value = "1.891.72"
$ws.Range("E12").Value = "  -4.28%  "

# Row 13
$ws.Range("D13").Value = "1.651.81"
$ws.Range("E13").Value = "  -4.58%  "

# Row 14
$ws.Range("E14").Value = "  -2.85%  "

# Row 15
$ws.Range("E15").Value = "  -0.96%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.04%  "

# Row 17
$ws.Range("D17").Value = "27.494.53"
$ws.Range("E17").Value = "  -2.84%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "239.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.54%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0729"
$ws.Range("E19").Value = "  -3.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.50%  "

# Row 21
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("E22").Value = "  -4.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.44%  "

# Row 24
$ws.Range("E24").Value = "  -2.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.50%  "

# Row 26
$ws.Range("E26").Value = "  -4.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "

# Row 29
$ws.Range("E29").Value = "  -2.54%  "

# Row 30
$ws.Range("E30").Value = "  -0.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0499"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.16%  "

# Row 33
$ws.Range("D33").Value = "1.442.22"
$ws.Range("E33").Value = "  -2.82%  "

# Row 34
$ws.Range("E34").Value = "  -5.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.11%  "

# Row 36
$ws.Range("E36").Value = "  -0.93%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.921"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0171"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.570"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.88%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.78%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.17%  "

# Row 43
$ws.Range("E43").Value = "  -3.93%  "

# Row 44
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.13%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.792"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "

# Row 46
$ws.Range("D46").Value = "1.800.13"
$ws.Range("E46").Value = "  -4.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.73%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.11%  "

# Row 49
$ws.Range("E49").Value = "  -6.84%  "

# Row 50
$ws.Range("E50").Value = "  -2.32%  "

# Row 51
$ws.Range("E51").Value = "  -5.28%  "
